$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country names that changed ranking order ---
$ws.Range("A55").Value = "Kirguistan"
$ws.Range("A56").Value = "Armenia"
$ws.Range("A123").Value = "Eslovaquia"
$ws.Range("A124").Value = "Surinam"
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Update case numbers (data refresh) ---
$ws.Range("B19").Value = 274525
$ws.Range("C19").Value = 2644
$ws.Range("D19").Value = 157635
$ws.Range("E19").Value = 113265
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 3625
$ws.Range("B26").Value = 137468
$ws.Range("C26").Value = 2345
$ws.Range("D26").Value = 91321
$ws.Range("E26").Value = 40076
$ws.Range("G26").Value = 50
$ws.Range("H26").Value = 6071
$ws.Range("B38").Value = 82924
$ws.Range("C38").Value = 181
$ws.Range("D38").Value = 77550
$ws.Range("E38").Value = 4812
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 562
$ws.Range("B48").Value = 55661
$ws.Range("C48").Value = 81
$ws.Range("E48").Value = 4585
$ws.Range("E53").Value = 3414
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 169
$ws.Range("B55").Value = 41645
$ws.Range("C55").Value = 272
$ws.Range("D55").Value = 33951
$ws.Range("E55").Value = 6201
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 1493
$ws.Range("B56").Value = 41495
$ws.Range("C56").Value = 196
$ws.Range("D56").Value = 34484
$ws.Range("E56").Value = 6194
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 817
$ws.Range("B88").Value = 9175
$ws.Range("C88").Value = 26
$ws.Range("D88").Value = 8831
$ws.Range("E88").Value = 219
$ws.Range("B111").Value = 4407
$ws.Range("C111").Value = 46
$ws.Range("D111").Value = 3487
$ws.Range("E111").Value = 853
$ws.Range("D122").Value = 2666
$ws.Range("E122").Value = 209
$ws.Range("B123").Value = 2855
$ws.Range("C123").Value = 54
$ws.Range("D123").Value = 1969
$ws.Range("E123").Value = 855
$ws.Range("H123").Value = 31
$ws.Range("B124").Value = 2838
$ws.Range("D124").Value = 1894
$ws.Range("E124").Value = 903
$ws.Range("H124").Value = 41
$ws.Range("B167").Value = 482
$ws.Range("C167").Value = 1
$ws.Range("E167").Value = 25
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# --- Update timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 11:55"
